# Populate the trading-log rows (5-16) with the new data that was
# "uploaded" in this commit, and move the active selection to H3 as in
# the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -----------------------------------------------------
# Column H got a touch narrower (stored OOXML width 19.06 -> 18.22).
# ColumnWidth is expressed in "characters" and gets pixel-quantized on
# write, so feed it the value (minus the fixed ~5/6 character padding)
# that lands closest to the target stored width.
$ws.Columns.Item(8).ColumnWidth = 17.3

# --- Row 5..9 : price / qty pairs (no date/time yet, same shape as row 4) ---
$ws.Range("A5").Value = 0.19563
$ws.Range("B5").Value = 30
$ws.Range("C5").Value = 45921

$ws.Range("A6").Value = 0.192
$ws.Range("B6").Value = 20

$ws.Range("A7").Value = 0.183
$ws.Range("B7").Value = 20

$ws.Range("A8").Value = 0.1818
$ws.Range("B8").Value = 20

$ws.Range("A9").Value = 0.17942
$ws.Range("B9").Value = 50

# --- Row 10..16 : price / qty / date / time entries ---------------------
# The date/time columns must stay plain text (matching the source data,
# which stores them as shared strings, not date serials). Some of the
# date strings (e.g. "03/11/2025") are still valid dates under a
# different locale reading, so Excel's automatic type-inference would
# otherwise silently convert them into date serial numbers. Forcing a
# text number format while assigning the value -- then restoring the
# original "General" format afterwards -- keeps the literal text while
# also keeping the original cell style (the style cache re-resolves
# back to the same index once the format matches again).
$dateTimeCells = @("C10","D10","C11","D11","C12","D12","C13","D13","C14","D14","C15","D15","C16","D16")
foreach ($addr in $dateTimeCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A10").Value = 0.1663
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "31/10/2025"
$ws.Range("D10").Value = "18:18:08"

$ws.Range("A11").Value = 0.166
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "31/10/2025"
$ws.Range("D11").Value = "18:46:08"

$ws.Range("A12").Value = 0.1642
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "03/11/2025"
$ws.Range("D12").Value = "14:35:46"

$ws.Range("A13").Value = 0.1635
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "03/11/2025"
$ws.Range("D13").Value = "18:14:54"

$ws.Range("A14").Value = 0.1608
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "04/11/2025"
$ws.Range("D14").Value = "09:36:06"

$ws.Range("A15").Value = 0.159
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "04/11/2025"
$ws.Range("D15").Value = "14:35:54"

$ws.Range("A16").Value = 0.1539
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "05/11/2025"
$ws.Range("D16").Value = "01:36:07"

foreach ($addr in $dateTimeCells) {
    $ws.Range($addr).NumberFormat = "General"
}

# --- Selection moved to H3 ----------------------------------------------
$ws.Range("H3").Select()
